$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = "paris"
$ws.Cells.Item(4, 2).Value = 4
$ws.Cells.Item(4, 3).Font.Bold = $false
$ws.Cells.Item(4, 4).Value = "DFT"
$ws.Cells.Item(4, 5).Value = "OTH"
$ws.Cells.Item(4, 6).Value = "53dcf950-aee9-43ba-bb93-9e7c5cd5833d"
$ws.Cells.Item(4, 7).Value = "By5SY2gA-_annotated.xlsx"
$ws.Cells.Item(4, 8).Value = "For instance, what about averaging WordNet path-based distance metrics and distance in word embedding space (for word similarity), and other ways of applying the affect data to email tone prediction?"
